$p = $ppt.ActivePresentation

# ------------------------------------------------------------------
# 1. Table on slide 16 (shape "Google Shape;213;p29") switches its
#    table style GUID from {3B3F232A-...} to {2548836E-...}.
# ------------------------------------------------------------------
$tableSlide = $p.Slides.Item(16)
$tableShape = $tableSlide.Shapes.Item(3)
$tableShape.Table.ApplyStyle("{2548836E-6285-4B6D-B24B-B31E2C39B1BA}")

# ------------------------------------------------------------------
# 2. Presentation theme color scheme swaps from "Integral" to the
#    stock "Office Theme" color values.
# ------------------------------------------------------------------
$theme = $p.SlideMaster.Theme
$colorScheme = $theme.ThemeColorScheme

$colorScheme.Item(1).RGB  = 0         # dk1      000000
$colorScheme.Item(2).RGB  = 16777215  # lt1      FFFFFF
$colorScheme.Item(3).RGB  = 6968388   # dk2      44546A
$colorScheme.Item(4).RGB  = 15132391  # lt2      E7E6E6
$colorScheme.Item(5).RGB  = 13998939  # accent1  5B9BD5
$colorScheme.Item(6).RGB  = 3243501   # accent2  ED7D31
$colorScheme.Item(7).RGB  = 10855845  # accent3  A5A5A5
$colorScheme.Item(8).RGB  = 49407     # accent4  FFC000
$colorScheme.Item(9).RGB  = 12874308  # accent5  4472C4
$colorScheme.Item(10).RGB = 4697456   # accent6  70AD47
$colorScheme.Item(11).RGB = 12673797  # hlink    0563C1
$colorScheme.Item(12).RGB = 7491477   # folHlink 954F72
